# Daily update - 12/15/2025 00:15:15
# Adds the new daily investment batch (date serial 46005 = 2025-12-15) to
# Daily_Profits_Calculations (Table1) and mirrors it into Daily_Report
# (Table4). Investor_Details' Total Profit Earned column recalculates
# automatically via its existing SUMIFS formulas.

$wb = $excel.ActiveWorkbook

$wsCalc   = $wb.Worksheets.Item("Daily_Profits_Calculations")
$wsReport = $wb.Worksheets.Item("Daily_Report")
$wsCharges = $wb.Worksheets.Item("Platofrm_Maintaince_Charges")

# ---------------------------------------------------------------------
# 1) New investor batch, dated 46005 (12/15/2025), total company pool 12000
# ---------------------------------------------------------------------
$rows = @(
    @{ Row=119; UserID="U001"; Name="Arun";                     Amt=2500 },
    @{ Row=120; UserID="U002"; Name="Grp Mem Mumtaj";           Amt=1500 },
    @{ Row=121; UserID="U003"; Name="Grp Mem Balram";           Amt=500  },
    @{ Row=122; UserID="U004"; Name="Bigil";                    Amt=1500 },
    @{ Row=123; UserID="U005"; Name="Vending Machine Referer";  Amt=1000 },
    @{ Row=124; UserID="U006"; Name="Grp Mem Athithya";         Amt=2000 },
    @{ Row=125; UserID="U007"; Name="Chinni";                   Amt=2000 },
    @{ Row=126; UserID="U008"; Name="Pradip";                   Amt=1000 }
)

$dateSerial = 46005
$companyTotal = 12000
$totalProfit = 434

foreach ($r in $rows) {
    $row = $r.Row

    $cellA = $wsCalc.Range("A$row")
    $cellA.Value = $dateSerial
    $cellA.NumberFormat = "m/d/yy"

    $wsCalc.Range("B$row").Value = $r.UserID
    $wsCalc.Range("C$row").Value = $r.Name
    $wsCalc.Range("G$row").Value = $r.Amt
    $wsCalc.Range("H$row").Value = $companyTotal
    $wsCalc.Range("I$row").Value = $totalProfit
    $wsCalc.Range("J$row").Formula = "=Table1[[#This Row],[User_Invest_Amount_As_On_Date]]/Table1[[#This Row],[Company_Total_Investment_As_On_Date]]"
    $wsCalc.Range("K$row").Formula = "=Table1[[#This Row],[Total_Profit]]*Table1[[#This Row],[Profit_%]]"
    $wsCalc.Range("L$row").Value = "Pending"
}

# ---------------------------------------------------------------------
# 2) Mirror the same batch into Daily_Report (rows 102-109). The Profit
#    figure is recomputed with the same arithmetic as Table1's J/K
#    columns (Total_Profit * (User_Invest_Amount/Company_Total_Invest))
#    so the cached values tie out exactly with Daily_Profits_Calculations.
# ---------------------------------------------------------------------
$reportRow = 102
foreach ($r in $rows) {
    $profitPct = $r.Amt / $companyTotal
    $profitValue = $totalProfit * $profitPct

    $cellA = $wsReport.Range("A$reportRow")
    $cellA.Value = $dateSerial
    $cellA.NumberFormat = "m/d/yy"

    $wsReport.Range("B$reportRow").Value = $r.UserID
    $wsReport.Range("C$reportRow").Value = $r.Amt
    $wsReport.Range("D$reportRow").Value = $companyTotal
    $wsReport.Range("E$reportRow").Value = $profitValue
    $wsReport.Range("F$reportRow").Value = "Pending"

    $reportRow++
}

# ---------------------------------------------------------------------
# 3) Small formatting clean-up on Daily_Report noticed in the same save:
#    a few legacy bordered/shaded cells get normalized to the plain
#    centered style already used by the rest of the table.
# ---------------------------------------------------------------------
$wsReport.Range("A5").Copy()
$wsReport.Range("A3").PasteSpecial(-4122)

$wsReport.Range("D5").Copy()
$wsReport.Range("D3").PasteSpecial(-4122)

$wsReport.Range("B5").Copy()
$wsReport.Range("F3").PasteSpecial(-4122)
$wsReport.Range("F5").PasteSpecial(-4122)
$wsReport.Range("F7").PasteSpecial(-4122)
$wsReport.Range("F9").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# 4) Restore view/selection state: Platofrm_Maintaince_Charges loses the
#    active tab, Daily_Profits_Calculations and Daily_Report get their
#    new selections, and Daily_Report becomes the active sheet.
# ---------------------------------------------------------------------
$wsCharges.Range("L11").Select()

$wsCalc.Activate()
$wsCalc.Range("K119:K126").Select()
$excel.ActiveWindow.ScrollRow = 106
$excel.ActiveWindow.ScrollColumn = 4

$wsReport.Activate()
$wsReport.Range("H102").Select()
$excel.ActiveWindow.ScrollRow = 82
$excel.ActiveWindow.ScrollColumn = 1
